# This workbook lists daily Fruit/Hortalizas prices for Mandarina at
# "Macroferia Regional de Talca". The commit adds two new daily price
# records (Murcott, Primera and Murcott, Segunda, dated 2023-10-13) at
# the top of the data block (rows 505-506), pushing the rest of the
# existing rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 505, shifting
# all the data that used to live at rows 505..558 down to 507..560.
$ws.Rows.Item(505).Insert()
$ws.Rows.Item(505).Insert()

# --- New row 505 -----------------------------------------------------
$ws.Range("A505").Value = 5
$ws.Range("B505").Value = "Macroferia Regional de Talca"
$ws.Range("C505").Value = "Maule"
$ws.Range("D505").Value = 45212
$ws.Range("E505").Value = 7
$ws.Range("F505").Value = "Fruta"
$ws.Range("G505").Value = 100102
$ws.Range("H505").Value = "Cítricos"
$ws.Range("I505").Value = 100102004
$ws.Range("J505").Value = "Mandarina"
$ws.Range("K505").Value = "Murcott"
$ws.Range("L505").Value = "Primera"
$ws.Range("M505").Value = 290
$ws.Range("N505").Value = 8000
$ws.Range("O505").Value = 8000
$ws.Range("P505").Value = 8000
$ws.Range("Q505").Value = "`$/bandeja 18 kilos"
$ws.Range("R505").Value = "Región de O'Higgins"
$ws.Range("S505").Value = 444
$ws.Range("T505").Value = 18

# --- New row 506 -----------------------------------------------------
$ws.Range("A506").Value = 5
$ws.Range("B506").Value = "Macroferia Regional de Talca"
$ws.Range("C506").Value = "Maule"
$ws.Range("D506").Value = 45212
$ws.Range("E506").Value = 7
$ws.Range("F506").Value = "Fruta"
$ws.Range("G506").Value = 100102
$ws.Range("H506").Value = "Cítricos"
$ws.Range("I506").Value = 100102004
$ws.Range("J506").Value = "Mandarina"
$ws.Range("K506").Value = "Murcott"
$ws.Range("L506").Value = "Segunda"
$ws.Range("M506").Value = 250
$ws.Range("N506").Value = 6000
$ws.Range("O506").Value = 6000
$ws.Range("P506").Value = 6000
$ws.Range("Q506").Value = "`$/bandeja 18 kilos"
$ws.Range("R506").Value = "Región de O'Higgins"
$ws.Range("S506").Value = 333
$ws.Range("T506").Value = 18

# Make sure the date cells use the same date/time number format as the
# rest of column D.
$ws.Range("D505:D506").NumberFormat = $ws.Range("D507").NumberFormat
